$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(30).Insert(-4121)

$ws.Cells.Item(30,1).Value = "CW3M"
$ws.Cells.Item(30,2).Value = "Baseline 2010-18_C151"
$ws.Cells.Item(30,3).Value = 2010
$ws.Cells.Item(30,19).Value = 2010
